$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# The GSC export window rolled forward by one day: the oldest date
# (2025-08-23, row 2) drops off and a new day (2025-11-21) is appended
# at the end. Deleting row 2 shifts every remaining row up by one,
# which is exactly the "Non-HTTPS/HTTPS URLs" shift seen in the diff,
# and it also removes the now-unused "2025-08-23" shared string.
$ws.Rows.Item(2).Delete()

# Append the new trailing day. Force the date to be written as literal
# text (matching how every other date cell in column A is stored),
# rather than letting Excel auto-convert the date-like string into a
# date serial number. A leading apostrophe is the standard Excel
# "treat as text" quote-prefix convention; ClearFormats then drops the
# quote-prefix cell style so the cell ends up styled like its neighbours.
$dateCell = $ws.Cells.Item(91, 1)
$dateCell.Value = "'2025-11-21"
$dateCell.ClearFormats()

$ws.Cells.Item(91, 2).Value = 0.0
$ws.Cells.Item(91, 3).Value = 0.0
